# Tilpasning til deling og oppdatering av halvårstall
# Fyll inn Mai (kolonne F) og Juni (kolonne G) punktlighetstall på "Total"-arket.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Total")

$ws.Range("F7").Value = 0.68
$ws.Range("G7").Value = 0.67

$ws.Range("F8").Value = 0.54
$ws.Range("G8").Value = 0.57

$ws.Range("F9").Value = 0.53
$ws.Range("G9").Value = 0.53

$ws.Range("F10").Value = 0.6
$ws.Range("G10").Value = 0.64

$ws.Range("F11").Value = 0.68
$ws.Range("G11").Value = 0.68

$ws.Range("F12").Value = 0.59
$ws.Range("G12").Value = 0.63

$ws.Range("F13").Value = 0.73
$ws.Range("G13").Value = 0.76

$ws.Range("F14").Value = 0.42
$ws.Range("G14").Value = 0.48

$ws.Range("F16").Value = 0.6
$ws.Range("G16").Value = 0.62

$ws.Range("F17").Value = 0.79
$ws.Range("G17").Value = 0.805

# Flytt markøren til siste celle som ble oppdatert, slik det fremstår i kilden.
$ws.Range("G17").Select() | Out-Null
